{"js": "// Locate the schedule table (first table in the document body).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Load existing rows so we can address them by index.\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// --- Step 1: simple text swaps on existing rows (buyer column = index 2) ---\n// Row indices include the header row at index 0.\n// 08:45 - 09:00 : ARMANDO VEL\u00c1SQUEZ -> COLFRESH COFFEE\nrows.items[2].getCell(2).value = \"COLFRESH COFFEE\";\n// 09:15 - 09:30 : INMERSSO BOUTIQUE -> FLOR A FRUTO\nrows.items[4].getCell(2).value = \"FLOR A FRUTO\";\n// 09:45 - 10:00 : COLFRESH COFFEE -> INMERSSO BOUTIQUE\nrows.items[6].getCell(2).value = \"INMERSSO BOUTIQUE\";\n// 10:00 - 10:15 : FLOR A FRUTO -> PROCOLOMBIA\nrows.items[7].getCell(2).value = \"PROCOLOMBIA\";\n\nawait context.sync();\n\n// --- Step 2: insert two new rows right after the 10:00 - 10:15 row (index 7) ---\nrows.items[7].insertRows(\"After\", 2, [\n  [\"10:15 - 10:30\", \"\", \"INTERLINK2AMERICAS\"],\n  [\"10:45 - 11:00\", \"\", \"ARMANDO VEL\u00c1SQUEZ\"],\n]);\nawait context.sync();\n\n// --- Step 3: remove the two now-trailing rows (old 11:15-11:30 and 12:00-12:15) ---\n// Re-load rows to get the fresh, post-insert collection.\nconst rows2 = table.rows;\nrows2.load(\"items\");\nawait context.sync();\n\nconst count = rows2.items.length;\n// Delete from the end first so indices of earlier rows stay valid.\nrows2.items[count - 1].delete();\nrows2.items[count - 2].delete();\nawait context.sync();\n", "ps1": "$doc = $word.ActiveDocument\n$table = $doc.Tables.Item(1)\n\n# --- Step 1: simple text swaps on existing rows (buyer column = column 3) ---\n# Row numbers are 1-based and include the header row (row 1).\n# 08:45 - 09:00 : ARMANDO VEL\u00c1SQUEZ -> COLFRESH COFFEE\n$table.Cell(3, 3).Range.Text = \"COLFRESH COFFEE\"\n# 09:15 - 09:30 : INMERSSO BOUTIQUE -> FLOR A FRUTO\n$table.Cell(5, 3).Range.Text = \"FLOR A FRUTO\"\n# 09:45 - 10:00 : COLFRESH COFFEE -> INMERSSO BOUTIQUE\n$table.Cell(7, 3).Range.Text = \"INMERSSO BOUTIQUE\"\n# 10:00 - 10:15 : FLOR A FRUTO -> PROCOLOMBIA\n$table.Cell(8, 3).Range.Text = \"PROCOLOMBIA\"\n\n# --- Step 2: insert two new rows right after the 10:00 - 10:15 row (row 8), i.e. before row 9 ---\n$refRow = $table.Rows.Item(9)\n$table.Rows.Add($refRow) | Out-Null\n$table.Rows.Add($refRow) | Out-Null\n\n$table.Cell(9, 1).Range.Text = \"10:15 - 10:30\"\n$table.Cell(9, 3).Range.Text = \"INTERLINK2AMERICAS\"\n$table.Cell(10, 1).Range.Text = \"10:45 - 11:00\"\n$table.Cell(10, 3).Range.Text = \"ARMANDO VEL\u00c1SQUEZ\"\n\n# --- Step 3: remove the two now-trailing rows (old 11:15-11:30 and 12:00-12:15) ---\n$table.Rows.Item($table.Rows.Count).Delete()\n$table.Rows.Item($table.Rows.Count).Delete()\n"}
